# Add two new columns (I: "کد درس", J: "کد ارائه") with header + per-row data
# for rows 2..122, matching the pattern:
#   I = 144001 + (row-2)
#   J = 133001 + (row-2)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 122

# 1) Copy the header cell's formatting (style index 1 -> red "Bad" style with border)
#    onto the whole new I:J range (header row + all data rows) before writing values,
#    so every new cell gets the same style used throughout the diff (s="1").
$ws.Range("A1").Copy()
$ws.Range("I1:J$lastRow").PasteSpecial(-4122)

# 2) Header row values (new shared strings "کد درس" / "کد ارائه")
$ws.Range("I1").Value = "کد درس"
$ws.Range("J1").Value = "کد ارائه"

# 3) Data rows: sequential numbering starting at 144001 / 133001 for row 2
for ($row = 2; $row -le $lastRow; $row++) {
    $ws.Cells.Item($row, 9).Value = 144000 + ($row - 1)
    $ws.Cells.Item($row, 10).Value = 133000 + ($row - 1)
}

# 4) New column widths for I (9) and J (10).
#    The engine rounds ColumnWidth to an internal 1/7-character pixel grid, so we pick
#    the input values whose rounded result lands exactly (I) / closest (J) to the
#    target stored widths of 13 and 12.75 respectively.
$ws.Columns.Item(9).ColumnWidth = 12.285714285714286
$ws.Columns.Item(10).ColumnWidth = 12.0

# 5) Update the visible selection to match the target (also clears the stale
#    scrolled "topLeftCell" position left over from the previous selection).
$ws.Application.Goto($ws.Range("A1"))
[void]$ws.Range("M8").Select()
